$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new Wins/Losses/Ties columns (AD/AE/AF), reusing
# the bold/centered/bordered formatting already used by the other header
# cells in row 1 (copy formatting from the adjacent AC1 header cell).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every data row (2-60) with
# the same W/L/T values.
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 30).Value = 58
    $ws.Cells.Item($r, 31).Value = 104
    $ws.Cells.Item($r, 32).Value = 0
}
